# Adapt tests to control version
# Add a "version" column to the settings sheet (header in C1, value 0 in C2)
# and make the settings sheet the active/selected sheet with C3 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")

# New "version" column next to form_title / form_id
$ws.Range("C1").Value = "version"
$ws.Range("C2").Value = 0

# Make "settings" the active sheet, selection on C3 (matches the tab/selection
# state captured in the saved workbook)
$ws.Activate()
$ws.Range("C3").Select()
